# Apply the post-testing parsing changes to the geo tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "plktest"
$ws.Range("B2").Value = "PLK1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "B. d."
$ws.Range("E2").Value = "dendrobatidis"
$ws.Range("F2").Value = "Swab"
$ws.Range("G2").Value = "Release"
$ws.Range("H2").Value = $true
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "plethodontidae"
$ws.Range("K2").Value = "Batrachoseps"
$ws.Range("L2").Value = "attenuatus"
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = "adult"
$ws.Range("O2").Value = 42326
$ws.Range("P2").Value = 37.878086000000003
$ws.Range("Q2").Value = -122.290059
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = "plk"

# --- Row 3 ---
$ws.Range("A3").Value = "plktest"
$ws.Range("B3").Value = "PLK2"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "B. d."
$ws.Range("E3").Value = "d"
$ws.Range("F3").Value = "Swab"
$ws.Range("G3").Value = "Release"
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = $false
$ws.Range("J3").Value = "plethodontidae"
$ws.Range("K3").Value = "Batrachoseps"
$ws.Range("L3").Value = "attenuatus"
$ws.Range("N3").Value = "adult"
$ws.Range("O3").Value = 42326
$ws.Range("P3").Value = 37.878086000000003
$ws.Range("Q3").Value = -122.290059
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = "plk"

# --- Row 4 is unchanged ---

# Move the active selection to I3, as left by the author after the edits.
$ws.Range("I3").Select() | Out-Null
